$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Apartment" name in A2 and A4/A5/A6 (shared-string table reorder
# causes displayed apartment codes to shift for these rows)
$ws.Range("A2").Value = "A11_02"
$ws.Range("A4").Value = "A10_04"
$ws.Range("A5").Value = "A09_01"
$ws.Range("A6").Value = "A10_01"

# Update "Su dung" (usage) values in column B
$ws.Range("B2").Value = 33
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 43
$ws.Range("B5").Value = 21
$ws.Range("B6").Value = 22
$ws.Range("B7").Value = 23
$ws.Range("B8").Value = 42
$ws.Range("B9").Value = 14
$ws.Range("B10").Value = 11

# Update selection to A2:A10 with active cell A2
$ws.Range("A2:A10").Select() | Out-Null
